$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column values are updated as plain text (matching the source data which stores
# prices as text, even when they look numeric). A leading apostrophe forces Excel to
# treat the entry as text instead of auto-converting it to a number/date, and resetting
# the cell style back to "Normal" afterwards clears the quote-prefix formatting flag
# that Excel applies, so the saved cell keeps no extra style attribute.

$d2 = $ws.Cells.Item(2, 4)
$d2.Value = "'68.395.13"
$d2.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -4.69%  "

$d3 = $ws.Cells.Item(3, 4)
$d3.Value = "'3.723.35"
$d3.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -3.83%  "

$d4 = $ws.Cells.Item(4, 4)
$d4.Value = "'0.998"
$d4.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.33%  "

$d5 = $ws.Cells.Item(5, 4)
$d5.Value = "'582.42"
$d5.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.50%  "

$d6 = $ws.Cells.Item(6, 4)
$d6.Value = "'183.43"
$d6.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +9.94%  "

$d7 = $ws.Cells.Item(7, 4)
$d7.Value = "'3.719.39"
$d7.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -3.66%  "

$d8 = $ws.Cells.Item(8, 4)
$d8.Value = "'0.629"
$d8.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -6.09%  "

$d9 = $ws.Cells.Item(9, 4)
$d9.Value = "'1.00"
$d9.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.00%  "

$d10 = $ws.Cells.Item(10, 4)
$d10.Value = "'0.715"
$d10.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -6.48%  "

$d11 = $ws.Cells.Item(11, 4)
$d11.Value = "'0.164"
$d11.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -8.63%  "

$d12 = $ws.Cells.Item(12, 4)
$d12.Value = "'53.99"
$d12.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.50%  "

$d13 = $ws.Cells.Item(13, 4)
$d13.Value = "'0.0000293"
$d13.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -8.93%  "

$d14 = $ws.Cells.Item(14, 4)
$d14.Value = "'10.53"
$d14.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -6.27%  "

$d15 = $ws.Cells.Item(15, 4)
$d15.Value = "'4.304.04"
$d15.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -4.59%  "

$d16 = $ws.Cells.Item(16, 4)
$d16.Value = "'3.715.44"
$d16.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -4.51%  "

$d17 = $ws.Cells.Item(17, 4)
$d17.Value = "'19.49"
$d17.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -7.03%  "

$ws.Cells.Item(18, 5).Value = "  -2.49%  "

$d19 = $ws.Cells.Item(19, 4)
$d19.Value = "'12.96"
$d19.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -6.43%  "

$d20 = $ws.Cells.Item(20, 4)
$d20.Value = "'1.13"
$d20.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -6.76%  "

$d21 = $ws.Cells.Item(21, 4)
$d21.Value = "'68.050.43"
$d21.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -5.10%  "

$d22 = $ws.Cells.Item(22, 4)
$d22.Value = "'408.74"
$d22.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -6.61%  "

$d23 = $ws.Cells.Item(23, 4)
$d23.Value = "'4.49"
$d23.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -3.90%  "

$d24 = $ws.Cells.Item(24, 4)
$d24.Value = "'88.57"
$d24.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -5.57%  "

$d25 = $ws.Cells.Item(25, 4)
$d25.Value = "'3.07"
$d25.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -5.92%  "

$d26 = $ws.Cells.Item(26, 4)
$d26.Value = "'12.86"
$d26.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -6.76%  "

$d27 = $ws.Cells.Item(27, 4)
$d27.Value = "'10.94"
$d27.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.41%  "

$d28 = $ws.Cells.Item(28, 4)
$d28.Value = "'3.83"
$d28.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -8.64%  "

$d29 = $ws.Cells.Item(29, 4)
$d29.Value = "'6.07"
$d29.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.50%  "

$d30 = $ws.Cells.Item(30, 4)
$d30.Value = "'9.58"
$d30.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -4.91%  "

$d31 = $ws.Cells.Item(31, 4)
$d31.Value = "'32.74"
$d31.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -6.37%  "

$d32 = $ws.Cells.Item(32, 4)
$d32.Value = "'7.66"
$d32.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.54%  "

$d33 = $ws.Cells.Item(33, 4)
$d33.Value = "'12.56"
$d33.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -7.08%  "

$ws.Cells.Item(34, 5).Value = "  -6.78%  "

$d35 = $ws.Cells.Item(35, 4)
$d35.Value = "'43.83"
$d35.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -11.35%  "

$d36 = $ws.Cells.Item(36, 4)
$d36.Value = "'608.34"
$d36.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.58%  "

$d37 = $ws.Cells.Item(37, 4)
$d37.Value = "'65.02"
$d37.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.69%  "

$d38 = $ws.Cells.Item(38, 4)
$d38.Value = "'0.0₃0903"
$d38.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -6.65%  "

$ws.Cells.Item(39, 5).Value = "  +0.20%  "

$d40 = $ws.Cells.Item(40, 4)
$d40.Value = "'0.401"
$d40.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -3.99%  "

$d41 = $ws.Cells.Item(41, 4)
$d41.Value = "'0.998"
$d41.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.51%  "

$ws.Cells.Item(42, 5).Value = "  -4.85%  "

$d43 = $ws.Cells.Item(43, 4)
$d43.Value = "'2.80"
$d43.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +8.22%  "

$d44 = $ws.Cells.Item(44, 4)
$d44.Value = "'3.03"
$d44.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -8.69%  "

$d45 = $ws.Cells.Item(45, 4)
$d45.Value = "'2.99"
$d45.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -5.00%  "

$d46 = $ws.Cells.Item(46, 4)
$d46.Value = "'0.0437"
$d46.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -6.57%  "

$d47 = $ws.Cells.Item(47, 4)
$d47.Value = "'9.26"
$d47.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -9.32%  "

$d48 = $ws.Cells.Item(48, 4)
$d48.Value = "'2.843.63"
$d48.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.97%  "

$ws.Cells.Item(49, 5).Value = "  -6.64%  "

$d50 = $ws.Cells.Item(50, 4)
$d50.Value = "'2.71"
$d50.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -4.20%  "

$d51 = $ws.Cells.Item(51, 4)
$d51.Value = "'3.07"
$d51.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -7.41%  "
